# TestData_Buildxact.xlsx edit:
#  - Users sheet: swap the PasswordBeforeReset/PasswordAfterReset values
#    (and their hyperlinks) between columns C and D for data rows 2 and 3,
#    then move the active-cell selection to B11.
#  - NewUsers sheet: move the active-cell selection to C9.

$wb = $excel.ActiveWorkbook
$wsUsers = $wb.Worksheets.Item("Users")
$wsNewUsers = $wb.Worksheets.Item("NewUsers")

# --- Users sheet: swap C2/D2 and C3/D3 cell values -------------------------
$c2 = $wsUsers.Range("C2").Value()
$d2 = $wsUsers.Range("D2").Value()
$wsUsers.Range("C2").Value = $d2
$wsUsers.Range("D2").Value = $c2

$c3 = $wsUsers.Range("C3").Value()
$d3 = $wsUsers.Range("D3").Value()
$wsUsers.Range("C3").Value = $d3
$wsUsers.Range("D3").Value = $c3

# --- Users sheet: hyperlinks follow their cell's new content ---------------
# (B2/D2 and C2 get swapped meanings too, matching the reset-password columns)
foreach ($hl in $wsUsers.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "`$B`$2") {
        $hl.Address = "mailto:Test@1234"
    } elseif ($addr -eq "`$D`$2") {
        $hl.Address = "mailto:shalini01@email.ghostinspector.com"
    } elseif ($addr -eq "`$C`$2") {
        $hl.Address = "mailto:Test@4321"
    } elseif ($addr -eq "`$C`$3") {
        $hl.Address = "mailto:Test@4321"
    } elseif ($addr -eq "`$D`$3") {
        $hl.Address = "mailto:Test@1234"
    }
}

# --- Selections (NewUsers first, Users last so "Users" stays the active tab)
$wsNewUsers.Range("C9").Select()
$wsUsers.Range("B11").Select()
